$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D5").Value = 2.471205793740484
$ws.Range("D6").Value = 0.06871335925147376
$ws.Range("D7").Value = -0.4023230641462986
$ws.Range("D8").Value = 0.09566701766772445
$ws.Range("D9").Value = 2.45035931757805
$ws.Range("D10").Value = 0.2911811583879726
$ws.Range("D11").Value = 2.475610245698677
$ws.Range("D12").Value = 0.389430763019782
$ws.Range("D13").Value = 0.3735002018570852
$ws.Range("D14").Value = 0.178700040371417
$ws.Range("D15").Value = 0.3641247559069888
$ws.Range("D16").Value = 0.06189974787140653
$ws.Range("D17").Value = -0.1046421997387022
$ws.Range("D18").Value = -0.01432397469862057
$ws.Range("D19").Value = 0.6225887028101049
$ws.Range("D20").Value = 0.3868634900551436
$ws.Range("D21").Value = 0.2096526231084588
$ws.Range("D22").Value = -0.02340715085628979
